# The underlying data rows (2-6) were reshuffled: each row's content moved
# to a different row position (a cyclic permutation), while the header row
# (row 1) and the columns that are identical across all of rows 2-6 stay
# untouched.
#
# Mapping of source row (current/before state) -> destination row (after
# state):
#   2 -> 5
#   3 -> 2
#   4 -> 3
#   5 -> 6
#   6 -> 4
#
# Only the columns whose values actually differ row-to-row need to move:
#   A (Id), B (Taxonsorteringsordning), E (TaxonId), F (Artnamn),
#   G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord),
#   AC (Publik kommentar)
# Restricting the copy to just these columns avoids disturbing other
# cells — in particular the date-like text cells in columns Y/AA, which
# would otherwise get silently reinterpreted as real Excel dates if
# rewritten through .Value2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 6; 6 = 4 }
$cols = @(1, 2, 5, 6, 7, 8, 17, 18, 29)   # A, B, E, F, G, H, Q, R, AC

# Capture the original values for the affected columns before overwriting
# anything, since sources and destinations overlap (it's a 5-cycle
# permutation). Using .Value2 (rather than .Formula/.Text) preserves full
# numeric precision for the coordinate columns (Q, R) instead of
# truncating to Excel's formula-bar ~15 significant-digit rendering.
$originalByRow = @{}
foreach ($srcRow in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Cells.Item($srcRow, $col).Value2
    }
    $originalByRow[$srcRow] = $rowValues
}

foreach ($srcRow in $mapping.Keys) {
    $dstRow = $mapping[$srcRow]
    $rowValues = $originalByRow[$srcRow]
    foreach ($col in $cols) {
        # Assigning a captured empty/blank Value2 back naturally clears the
        # destination cell (e.g. AC for rows whose source had no comment),
        # so no extra ClearContents() step is needed.
        $ws.Cells.Item($dstRow, $col).Value2 = $rowValues[$col]
    }
}
